# Update countries & provincias Spain
# Applies the 21-Sep-2020 10:05 data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp in A1
#  - updates case counters for several countries (Rusia, Filipinas, Singapur,
#    Armenia, Bulgaria/Hungria, Aruba/Mayotte/Gambia/Tailandia/Georgia,
#    Estonia, Letonia)
#  - re-sorts the table (descending by total cases), which swaps the
#    row order of a few countries whose totals crossed over

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes)
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- "last updated" banner ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 10:05"

# --- updated country rows ---
Set-Row 7   "Rusia"          1109595 6196 911973 178133 0 71 19489
Set-Row 24  "Filipinas"      290190  3475 230233 54958  0 15 4999
Set-Row 57  "Singapur"       57607   31   57181  399    0 0  27
Set-Row 63  "Armenia"        47552   121  42637  3979   0 6  936

# Hungria overtakes Bulgaria
Set-Row 84  "Hungria"        18866   876  4401   13779  0 3  686
Set-Row 85  "Bulgaria"       18863   0    13580  4522   0 0  761

# Georgia overtakes Aruba, Mayotte, Gambia, Tailandia
Set-Row 133 "Georgia"        3695    193  1534   2141   0 1  20
Set-Row 134 "Aruba"          3551    0    2239   1289   0 0  23
Set-Row 135 "Mayotte"        3541    0    2964   537    0 0  40
Set-Row 136 "Gambia"         3526    0    1992   1426   0 0  108
Set-Row 137 "Tailandia"      3506    0    3342   105    0 0  59

Set-Row 144 "Estonia"        2941    17   2379   498    0 0  64
Set-Row 161 "Letonia"        1526    1    1248   242    0 0  36

# Timor Oriental / Santa Lucia swap order (tie on total cases)
Set-Row 204 "Timor Oriental" 27      0    26     1      0 0  0
Set-Row 205 "Santa Lucia"    27      0    26     1      0 0  0

# Islas Malvinas / Montserrat swap order
Set-Row 214 "Islas Malvinas" 13      0    13     0      0 0  0
Set-Row 215 "Montserrat"     13      0    12     0      0 0  1
